# MergeTables.py gemaakt die Conditions_table.xlsx met experiment_table.csv fuseert
# Replays (via COM) the effect of merging in the experiment_table.csv data:
#  - the condition rows get re-sorted by Wells (column F), descending
#  - a handful of Titr0/Titr2 Cmin/Cmax cells receive updated, merged values
#    and a centered number format
#  - sheet/window view + page setup cosmetics

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Sort the data block A2:R12 by column F (Wells), descending ---------
$sortObj = $ws.Sort
$sortObj.SortFields.Clear()
$sortObj.SortFields.Add($ws.Range("F2:F12"), $null, 2, $null, 0)
$sortObj.SetRange($ws.Range("A1:R12"))
$sortObj.Header = 1
$sortObj.Apply()

# --- 2. Overwrite the cells whose values came from the merged CSV ----------
# row 3  = LapaBiniVino120 (Titr0 Cmin/Cmax)
$ws.Range("P3").Value = 0.11903712
$ws.Range("Q3").Value = 0.11903712
# row 7  = LapaVino120 (Titr2 Cmin/Cmax)
$ws.Range("L7").Value = 0.11903712
$ws.Range("M7").Value = 0.11903712
# row 9  = BiniVino120 (Titr2 Cmin/Cmax)
$ws.Range("L9").Value = 0.11903712
$ws.Range("M9").Value = 0.11903712
# row 11 = Pos_ctrl (Titr1 Cmin/Cmax)
$ws.Range("H11").Value = 10.0189576
$ws.Range("I11").Value = 10.0189576
# row 12 = Vino120 (Titr1 Cmin/Cmax)
$ws.Range("H12").Value = 0.119037158575
$ws.Range("I12").Value = 0.119037158575

# Those ten cells also pick up a centered-alignment style
foreach ($addr in @("P3","Q3","L7","M7","L9","M9","H11","I11","H12","I12")) {
    $ws.Range($addr).HorizontalAlignment = -4108
}

# --- 3. View / selection cosmetics -----------------------------------------
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A1:R12").Select() | Out-Null

# --- 4. Page setup -----------------------------------------------------------
$ps = $ws.PageSetup
$ps.PaperSize = 9
$ps.Orientation = 1
